# Prepare public release and harden processing reliability
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Sheet1" to "Sheet"
$ws.Name = "Sheet"

# A2: "Item1" -> "Test"
$ws.Range("A2").Value = "Test"

# B2: numeric 100 -> text "123".
# A leading apostrophe forces Excel to store the value as text instead of
# coercing the numeric-looking string back into a number; reset the style
# afterwards so the cell doesn't keep the "quote prefix" formatting flag.
$ws.Range("B2").Value = "'123"
$ws.Range("B2").Style = "Normal"

# Drop row 3 entirely ("Item2" / 200), shrinking the used range to A1:B2
$ws.Range("A3:B3").EntireRow.Delete()
